$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.252.72"
$ws.Range("E2").Value = "  +1.59%  "
$ws.Range("D3").Value = "1.907.69"
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "307.50"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.80%  "
$ws.Range("E6").Value = "  +0.02%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5256"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +3.58%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3784"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +3.64%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07261"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.55%  "
$ws.Range("E10").Value = "  +3.66%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.9004"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.31%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08269"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +10.47%  "
$ws.Range("D13").Value = "1.914.21"
$ws.Range("E13").Value = "  +2.78%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "95.33"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.48%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.280"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.30%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.001"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008605"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.54%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "14.50"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.63%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.9999"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.01%  "
$ws.Range("D20").Value = "27.283.32"
$ws.Range("E20").Value = "  +1.60%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.065"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.56%  "
$ws.Range("D22").Value = "2.153.05"
$ws.Range("E22").Value = "  +2.05%  "
$ws.Range("E23").Value = "  +3.20%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.457"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.81%  "
$ws.Range("E25").Value = "  +10.44%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "145.99"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.744"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.58%  "
$ws.Range("E28").Value = "  +1.74%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "114.80"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.08%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.970"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +5.89%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.812"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.99%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09201"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.89%  "
$ws.Range("B33").Value = "Hedera"
$ws.Range("C33").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05087"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.22%  "
$ws.Range("B34").Value = "ImmutableX"
$ws.Range("C34").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.8031"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +7.57%  "
$ws.Range("E35").Value = "  +8.09%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.949"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.24%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.361"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +4.75%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.569"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.03%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.5748"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.78%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01978"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.075"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.36%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "9.053"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +5.79%  "
$ws.Range("E43").Value = "  +1.29%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "118.54"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.30%  "
$ws.Range("E45").Value = "  +2.33%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4841"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.24%  "
$ws.Range("E47").Value = "  +0.05%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "10.16"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.47%  "
$ws.Range("E49").Value = "  +4.27%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "37.63"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.66%  "
$ws.Range("E51").Value = "  +1.56%  "
